# "added some p values"
#
# The author filled in p-value columns E:J (A->G Syn CpG-v-NonCpG,
# A->G NonSyn CpG-v-NonCpG, A->G Syn-v-NonSyn, T->C Syn CpG-v-NonCpG,
# T->C NonSyn CpG-v-NonCpG, T->C Syn-v-NonSyn) for several data rows.
# Wherever a test wasn't significant/applicable the author typed the
# literal text "<0.01"; otherwise an actual numeric p-value was entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Dengue 1
$ws.Range("E2").Value = "<0.01"
$ws.Range("F2").Value = "<0.01"
$ws.Range("G2").Value = "<0.01"
$ws.Range("H2").Value = "<0.01"
$ws.Range("I2").Value = "<0.01"
$ws.Range("J2").Value = "<0.01"

# Row 3 - Dengue 2
$ws.Range("E3").Value = "<0.01"
$ws.Range("F3").Value = "<0.01"
$ws.Range("G3").Value = "<0.01"
$ws.Range("H3").Value = "<0.01"
$ws.Range("I3").Value = 0.061
$ws.Range("J3").Value = "<0.01"

# Row 4 - Dengue 3
$ws.Range("E4").Value = "<0.01"
$ws.Range("F4").Value = 0.0227
$ws.Range("G4").Value = "<0.01"
$ws.Range("H4").Value = "<0.01"
$ws.Range("I4").Value = "<0.01"
$ws.Range("J4").Value = "<0.01"

# Row 5 - Dengue 4
$ws.Range("E5").Value = "<0.01"
$ws.Range("F5").Value = 0.073
$ws.Range("G5").Value = "<0.01"
$ws.Range("H5").Value = "<0.01"
$ws.Range("I5").Value = "<0.01"
$ws.Range("J5").Value = "<0.01"

# Row 6 - Bk Polyoma VP1
$ws.Range("E6").Value = 0.0216
$ws.Range("F6").Value = 0.382
$ws.Range("G6").Value = "<0.01"
$ws.Range("H6").Value = 0.0237
$ws.Range("I6").Value = 0.198
$ws.Range("J6").Value = "<0.01"

# Row 7 - Hepatitis B PTP
$ws.Range("E7").Value = "<0.01"
$ws.Range("F7").Value = 0.883
$ws.Range("G7").Value = "<0.01"
$ws.Range("H7").Value = "<0.01"
$ws.Range("I7").Value = "<0.01"
$ws.Range("J7").Value = "<0.01"

# Row 8 - Rhino C
$ws.Range("E8").Value = "<0.01"
$ws.Range("F8").Value = "<0.01"
$ws.Range("G8").Value = "<0.01"
$ws.Range("H8").Value = "<0.01"
$ws.Range("I8").Value = 0.163
$ws.Range("J8").Value = "<0.01"

# Row 9 - Human Parainfluenza 1 F
$ws.Range("E9").Value = "<0.01"
$ws.Range("F9").Value = "<0.01"
$ws.Range("G9").Value = 0.988
$ws.Range("H9").Value = 0.0559
$ws.Range("I9").Value = 0.159
$ws.Range("J9").Value = 1

# Row 11 - Human Parainfluenza 3 HN
$ws.Range("E11").Value = 0.103
$ws.Range("F11").Value = 0.157
$ws.Range("G11").Value = "<0.01"
$ws.Range("H11").Value = 0.0772
$ws.Range("I11").Value = 0.139
$ws.Range("J11").Value = "<0.01"

# Row 14 - Influenza A NA H3N2
$ws.Range("E14").Value = "<0.01"
$ws.Range("F14").Value = 0.224
$ws.Range("G14").Value = "<0.01"
$ws.Range("H14").Value = "<0.01"
$ws.Range("I14").Value = 0.126
$ws.Range("J14").Value = "<0.01"

# Row 15 - Influenza A HA H1N1
$ws.Range("E15").Value = "<0.01"
$ws.Range("F15").Value = 0.385
$ws.Range("G15").Value = "<0.01"
$ws.Range("H15").Value = "<0.01"
$ws.Range("I15").Value = "<0.01"
$ws.Range("J15").Value = "<0.01"

# Row 16 - Influenza A HA H3N2
$ws.Range("E16").Value = "<0.01"
$ws.Range("F16").Value = 0.0677
$ws.Range("G16").Value = "<0.01"
$ws.Range("H16").Value = "<0.01"
$ws.Range("I16").Value = "<0.01"
$ws.Range("J16").Value = "<0.01"

# Row 17 - Influenza A NA H1N1
$ws.Range("E17").Value = "<0.01"
$ws.Range("F17").Value = 0.19
$ws.Range("G17").Value = "<0.01"
$ws.Range("H17").Value = "<0.01"
$ws.Range("I17").Value = "<0.01"
$ws.Range("J17").Value = "<0.01"

# Row 18 - Influenza B HA
$ws.Range("E18").Value = "<0.01"
$ws.Range("F18").Value = 0.0946
$ws.Range("G18").Value = "<0.01"
$ws.Range("H18").Value = "<0.01"
$ws.Range("I18").Value = 0.162
$ws.Range("J18").Value = "<0.01"

# Row 19 - Influenza B NA
$ws.Range("E19").Value = "<0.01"
$ws.Range("F19").Value = 0.0328
$ws.Range("G19").Value = "<0.01"
$ws.Range("H19").Value = "<0.01"
$ws.Range("I19").Value = 0.237
$ws.Range("J19").Value = "<0.01"

# Row 20 - Entero D68 VP1
$ws.Range("E20").Value = "<0.01"
$ws.Range("F20").Value = 0.37
$ws.Range("G20").Value = "<0.01"
$ws.Range("H20").Value = 0.0219
$ws.Range("I20").Value = 0.196
$ws.Range("J20").Value = "<0.01"

# Match the author's final cursor position/selection (cell J20).
$ws.Range("J20").Select()

# Slightly narrow column A, matching the author's manual resize.
$ws.Columns.Item(1).ColumnWidth = 22.75
